# Apply committee-votes data corrections (sponsorship/committee validation pass)
# as described in the commit: "validation and updates to the core forge model"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 0

# Row 4
$ws.Range("V4").Value = 0

# Row 6
$ws.Range("B6").Value = 0
$ws.Range("Q6").Value = 0

# Row 7
$ws.Range("J7").Value = 1

# Row 8
$ws.Range("J8").Value = 0
$ws.Range("O8").Value = 0.6666666666666666
$ws.Range("P8").Value = 0.8
$ws.Range("U8").Value = 0.8333333333333334
$ws.Range("AA8").Value = 0.8333333333333334
$ws.Range("AB8").Value = 0.8

# Row 9
$ws.Range("W9").Value = 0.625
$ws.Range("Z9").Value = 0.7142857142857143
$ws.Range("AD9").Value = 0.8571428571428571

# Row 10
$ws.Range("G10").Value = 1
$ws.Range("H10").Value = 0
$ws.Range("O10").Value = 0
$ws.Range("Z10").Value = 0

# Row 11
$ws.Range("O11").Value = 0.5
$ws.Range("P11").Value = 0.6
$ws.Range("U11").Value = 0.8333333333333334

# Row 15
$ws.Range("H15").Value = 0.6666666666666666
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0.5
$ws.Range("P15").Value = 0.6666666666666666
$ws.Range("R15").Value = 0
$ws.Range("U15").Value = 1
$ws.Range("AA15").Value = 1
$ws.Range("AB15").Value = 0.6666666666666666

# Row 16
$ws.Range("H16").Value = 0.8
$ws.Range("K16").Value = 0.6
$ws.Range("O16").Value = 0.6666666666666666
$ws.Range("R16").Value = 0.6666666666666666
$ws.Range("U16").Value = 0.6666666666666666
$ws.Range("AA16").Value = 0.6
$ws.Range("AB16").Value = 0.4

# Row 17
$ws.Range("F17").Value = 0

# Row 18
$ws.Range("O18").Value = 0
$ws.Range("P18").Value = 0.6666666666666666
$ws.Range("U18").Value = 0.75
$ws.Range("AA18").Value = 0.75
$ws.Range("AB18").Value = 0.75

# Row 21
$ws.Range("H21").Value = 0.8333333333333334
$ws.Range("K21").Value = 0.8333333333333334
$ws.Range("O21").Value = 1
$ws.Range("P21").Value = 0.6666666666666666
$ws.Range("R21").Value = 0.75
$ws.Range("AB21").Value = 0.8888888888888888

# Row 22
$ws.Range("D22").Value = 0
$ws.Range("Z22").Value = 0

# Row 23
$ws.Range("I23").Value = 0.625
$ws.Range("Z23").Value = 0.75
$ws.Range("AD23").Value = 0.7142857142857143

# Row 26
$ws.Range("I26").Value = 0.7142857142857143
$ws.Range("J26").Value = 0
$ws.Range("V26").Value = 0
$ws.Range("W26").Value = 0.75
$ws.Range("AD26").Value = 0.6666666666666666

# Row 27
$ws.Range("H27").Value = 0.8333333333333334
$ws.Range("O27").Value = 1
$ws.Range("P27").Value = 0.6
$ws.Range("R27").Value = 0.75

# Row 28
$ws.Range("H28").Value = 0.8
$ws.Range("O28").Value = 0.6666666666666666
$ws.Range("P28").Value = 0.4
$ws.Range("R28").Value = 0.75
$ws.Range("U28").Value = 0.8888888888888888

# Row 30
$ws.Range("I30").Value = 0.8571428571428571
$ws.Range("W30").Value = 0.7142857142857143
$ws.Range("Z30").Value = 0.6666666666666666
